# Rearranges the data rows (3-14) of the "Artfynd" sheet: the full content
# of each row is replaced by the full content of another row, following a
# fixed permutation (three independent row-content rotation cycles). Row 1
# (header) and row 2 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row number -> source row number (content of source row moves into target row)
$map = @{
    3  = 9
    4  = 3
    5  = 10
    6  = 11
    7  = 4
    8  = 12
    9  = 13
    10 = 5
    11 = 14
    12 = 6
    13 = 7
    14 = 8
}

$firstRow = 3
$lastRow = 14
$lastCol = 51   # column AY

$srcRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$origVals = $srcRange.Value()

$numRows = $lastRow - $firstRow + 1
# NOTE: Range.Value() returns a 1-based array ([1..numRows, 1..lastCol]); a
# freshly `New-Object`-ed array is a normal 0-based .NET array. Excel's
# Range.Value setter only accepts an array whose size exactly matches the
# target range's dimensions (it does NOT care about the lower bound), so we
# build a 0-based array of the same size and just offset our own indices by
# one when reading from $origVals.
$newVals = New-Object 'object[,]' $numRows, $lastCol

for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $map[$targetRow]
    $targetIdx = $targetRow - $firstRow
    $sourceIdx = $sourceRow - $firstRow + 1
    for ($col = 1; $col -le $lastCol; $col++) {
        $newVals[$targetIdx, ($col - 1)] = $origVals[$sourceIdx, $col]
    }
}

$srcRange.Value = $newVals
